$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F ("dSF") values for the specified rows per repull/push of data
$updates = @{
    2  = 3
    10 = 0
    13 = 0
    17 = 2
    19 = -4
    20 = -1
    26 = -3
    27 = 3
    29 = -4
    30 = -1
    34 = 1
    36 = -2
    37 = 3
    41 = 5
    42 = 1
    46 = 5
    49 = -3
    50 = -6
    52 = 7
    54 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
